# Daily Satellite Data Update
# Refresh the pass-prediction times and the "Oblačnost" (cloud-cover)
# bucket counts/colors for the three dates on the "Přelety" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (25.12.2025) -----------------------------------------------
$ws.Range("C2").Value = "00:37"
$ws.Range("E2").Value = "05:48:36"
$ws.Range("F2").Value = "05:52:49"
$ws.Range("G2").Value = "05:53:07"
$ws.Range("H2").Value = "05:53:26"
$ws.Range("I2").Value = "05:57:38"

# --- Row 3 (27.12.2025) -------------------------------------------------
$ws.Range("D3").Value = "10:43"
$ws.Range("E3").Value = "06:15:33"
$ws.Range("F3").Value = "06:18:00"
$ws.Range("G3").Value = "06:20:54"
$ws.Range("H3").Value = "06:23:49"
$ws.Range("I3").Value = "06:26:16"

# --- Row 4 (28.12.2025) -------------------------------------------------
$ws.Range("C4").Value = "05:36"
$ws.Range("D4").Value = "10:36"
$ws.Range("E4").Value = "05:40:42"
$ws.Range("F4").Value = "05:43:12"
$ws.Range("G4").Value = "05:45:59"
$ws.Range("H4").Value = "05:48:48"
$ws.Range("I4").Value = "05:51:18"

# --- Oblačnost (cloud coverage) counters ---------------------------------
# The "Střední" (medium) bucket fill is refreshed to a slightly darker blue.
$ws.Range("Q2").Interior.Color = 13341278   # RGB(0x5E,0x92,0xCB) -> BGR int

# Row 2 counts
$ws.Range("O2").Value = 4
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 0

# Row 3: now matches the "Nízká" (low) / "Střední" (medium) look used on
# row 2, so pull the already-updated formatting across before setting the
# new counts.
$ws.Range("O2").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("Q2").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 3

# Row 4: counts are unchanged; only P4's style index shifts downstream once
# the unused "B2CBE6" fill/style drops out, so no direct edit is needed
# here beyond what already matches.
